# Fix "Recorded By" (column G) entries so that "System"/"system" is no
# longer listed first in the comma-separated list of recorders; instead
# it is moved to the end of the list (other entries keep their relative
# order). Rows whose list does not include "system" (case-insensitive),
# or that only contain a single entry, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($i = 2; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, 7)   # Column G = "Recorded By"
    $v = $cell.Value2

    if ($v -ne $null -and [string]$v -ne "") {
        $parts = [string]$v -split ", "

        if ($parts.Count -gt 1) {
            $hasSystem = $false
            foreach ($p in $parts) {
                if ($p.Trim().ToLower() -eq "system") {
                    $hasSystem = $true
                }
            }

            if ($hasSystem) {
                $reversed = $parts[($parts.Count - 1)..0]
                $newVal = $reversed -join ", "
                $cell.Value2 = $newVal
            }
        }
    }
}
